$wb = $excel.ActiveWorkbook

# --- Overview sheet: update summary row for b.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 16:56:33"

# --- zh-cn sheet: a new handoff was generated for b.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Copy F2 (already a literal text "False") into F3 so Excel doesn't
# auto-convert the string into a native Boolean value.
$zhcn.Range("F2").Copy($zhcn.Range("F3"))
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-07 16:56:27"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bc43aabda39d2bade88bab34569e4e1ed1ac09e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fe3abc2da7243d4a324df82a072f918fd4ea879/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: a new handoff was generated for b.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F2").Copy($dede.Range("F3"))
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-07 16:56:33"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bc43aabda39d2bade88bab34569e4e1ed1ac09e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fe3abc2da7243d4a324df82a072f918fd4ea879/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
